# daily auto push: 2026-02-22 22:41 UTC
# Insert a new data row (2026/02/23, 月, 5, 200) right after the existing
# "2026/02/23" row (row 864), pushing the remaining rows down by one and
# extending the used range to D906.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 864; everything from 864..905 shifts to 865..906.
$ws.Rows.Item(864).Insert()

# Write the new row's values. The date column is plain text in this sheet
# (e.g. "2026/02/23"), so a direct .Value assignment would get
# auto-converted into a date serial by Excel. Instead, write it as a
# formula that evaluates to the literal text, then paste-special just the
# value back over itself so the cell ends up holding a plain text value
# (no formula, no stray number-format/style changes).
$ws.Range("A864").Formula = "=""2026/02/23"""
$ws.Range("A864").Copy() | Out-Null
$ws.Range("A864").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("B864").Value = "月"
$ws.Range("C864").Value = 5
$ws.Range("D864").Value = 200
